$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: datetime_utc_matlab
$ws.Range("A2").Value = "datetime_utc_matlab"
$ws.Range("B2").Value = "PI-provided UTC date and time "
$ws.Range("C2").Value = "Date"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "YYYY-MM-DD hh:mm:ss"

# Row 3: latitude_matlab
$ws.Range("A3").Value = "latitude_matlab"
$ws.Range("B3").Value = "Latitude of sample event provided by PI"
$ws.Range("C3").Value = "numeric"
$ws.Range("D3").Value = "degree"

# Row 4: longitude_matlab
$ws.Range("A4").Value = "longitude_matlab"
$ws.Range("B4").Value = "Longitude of sample event provided by PI"
$ws.Range("C4").Value = "numeric"
$ws.Range("D4").Value = "degree"

# Row 5: latitude_API
$ws.Range("A5").Value = "latitude_API"
$ws.Range("B5").Value = "Latitude of sample event provided by NES-LTER API"
$ws.Range("C5").Value = "numeric"
$ws.Range("D5").Value = "degree"

# Row 6: longitude_API
$ws.Range("A6").Value = "longitude_API"
$ws.Range("B6").Value = "Longitude of sample event provided by NES-LTER API"
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "degree"

# Row 7: depth
$ws.Range("A7").Value = "depth"
$ws.Range("B7").Value = "Depth of sample below sea surface. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$ws.Range("C7").Value = "numeric"
$ws.Range("D7").Value = "meter"
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""

# Row 8: biosat
$ws.Range("A8").Value = "biosat"
$ws.Range("B8").Value = "Percent biological saturation, the oxygen-argon ratio divided by the equilibrium value of that ratio "
$ws.Range("C8").Value = "numeric"
$ws.Range("D8").Value = "dimensionless"
$ws.Range("F8").Value = "NaN"
$ws.Range("G8").Value = "Missing value"

# Row 9: O2_Ar_ratio (new row)
$ws.Range("A9").Value = "O2_Ar_ratio"
$ws.Range("B9").Value = "Oxygen-argon ratio of EIMS sample from underway corrected for air values"
$ws.Range("C9").Value = "numeric"
$ws.Range("D9").Value = "dimensionless"
$ws.Range("F9").Value = "NaN"
$ws.Range("G9").Value = "Missing value"

# Column A width adjustment (target stored width 24.3984375 chars; engine quantizes
# column widths to 1/6-character steps, so 23.5 is the closest achievable input and
# lands on the nearest representable stored width, 24.333333333333332)
$ws.Columns("A").ColumnWidth = 23.5

# View changes: scroll so column C is the leftmost visible column, and select E2
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E2").Select()
